$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Remove the three obsolete data rows (id 1,2,3 -> MACAULAY CULKIN, ALFONSO RIBEIRO,
#    DAVID HASSELHOFF). The two remaining rows (id 4, id 5) shift up to rows 4 and 5.
$ws.Range("A4:D6").EntireRow.Delete()

# 2. Remove the "Documento" and "correo" columns (C and D) entirely.
$ws.Range("C1:D1").EntireColumn.Delete()

# 3. Update the report title and the first column header.
$ws.Range("A1").Value = "Informe de Actividad 16-11-2018"
$ws.Range("A3").Value = "matricula"

# 4. Re-apply the AutoFilter over the new, narrower header range (A3:B3).
#    AutoFilter() toggles the filter on/off, so only call it once if it is
#    currently off, to land on an "on" state that references A3:B3.
if ($ws.AutoFilterMode()) {
    $ws.Range("A3:B3").AutoFilter()
}
$ws.Range("A3:B3").AutoFilter()

# 5. Fix up the hidden _FilterDatabase defined name so it also points at A3:B3.
$fd = $wb.Names.Item(1)
$fd.RefersTo = "='Mi hoja guay'!`$A`$3:`$B`$3"

# 6. Best-fit columns A and B to their (now shorter) content and move the
#    selection to the last populated cell, mirroring the original session.
$ws.Columns.Item(1).EntireColumn.AutoFit()
$ws.Columns.Item(2).EntireColumn.AutoFit()
$ws.Columns.Item(1).ColumnWidth = 9.1666666666667
$ws.Columns.Item(2).ColumnWidth = 13.1666666666667
$ws.Range("B5").Select()

# 7. Keep gridlines/headers visible, matching the source sheet view flags.
$win = $excel.ActiveWindow
$win.DisplayGridlines = $true
$win.DisplayHeadings = $true
